$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 / Row 14: WrappedEther and Polkadot swap rank positions (name + link),
# with freshly updated price/volume figures for each.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.528.76"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.03"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.34"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("E6").Value = "  +2.12%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.250"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "  -0.63%  "

$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.82"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.18"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("B13").Value = "WrappedEther"

$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.630.13"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("B14").Value = "Polkadot"

$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.15"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.28"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = "  +3.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.556.94"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.47"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  +3.22%  "

$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  +13.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.48"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.92"
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("D30").ClearFormats()

$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = "  +4.19%  "

$ws.Range("E33").Value = "  +1.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.253.66"
$ws.Range("D34").ClearFormats()

$ws.Range("E34").Value = "  +7.40%  "

$ws.Range("E35").Value = "  +0.54%  "

$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("E37").Value = "  +4.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.511"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = "  +1.54%  "

$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.797"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = "  -1.27%  "

$ws.Range("E41").Value = "  -2.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.798"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.764.41"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.26"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("E46").Value = "  +3.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.04"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0960"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = "  +2.30%  "

$ws.Range("E51").Value = "  -0.52%  "
